$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new draw result row (row 22), matching the existing rows which
# store every column (including numeric-looking ones) as text.
$newRow = $ws.Range("A22:E22")

# Force text storage first so the numeric-looking strings ("251008") are
# not auto-converted into real numbers by Excel's type inference.
$newRow.NumberFormat = "@"

$ws.Range("A22").Value = "2025-10-08"
$ws.Range("B22").Value = "Pick 3"
$ws.Range("C22").Value = "251008"
$ws.Range("D22").Value = "4-7-0"
$ws.Range("E22").Value = "2025-10-08T21:38:15.296+04:00"

# Reset the cell style back to Normal (style 0), same as every other row
# in the sheet, so only the value/type changed - no stray text-format
# style index is introduced.
$newRow.Style = "Normal"
